$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 7: columns G (7) through BB (54) - odds updated for Blooming - Aurora
$row7 = @(2.3, 3.1, 3.3, 3.1, 2, 4, 1.08, 8, 1.36, 3, 2.25, 1.62, 1.5, 2.5, 1.91, 1.8, 7, 10, 10, 21, 21, 34, 7.5, 6, 15, 51, 351, 8.5, 15, 12, 34, 29, 41, 4.33, 13, 26, 41, 81, 201, 2.5, 8.5, 67, 5, 19, 29, 67, 101, 251)
for ($i = 0; $i -lt $row7.Length; $i++) {
    $ws.Cells.Item(7, 7 + $i).Value = $row7[$i]
}

# Row 8: columns G (7) through BB (54) - odds updated for Corinthians - Bahia
$row8 = @(1.67, 4, 4.75, 2.25, 2.38, 4.75, 1.04, 13, 1.22, 4.33, 1.67, 2.15, 1.33, 3.25, 1.7, 2.05, 8.5, 9, 8.5, 13, 13, 21, 15, 8, 15, 41, 151, 15, 26, 15, 51, 34, 34, 3.75, 8.5, 17, 26, 41, 101, 3.25, 7.5, 51, 6.5, 23, 29, 81, 81, 151)
for ($i = 0; $i -lt $row8.Length; $i++) {
    $ws.Cells.Item(8, 7 + $i).Value = $row8[$i]
}
